$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 648, shifting existing rows 648:755 down to 649:756
$ws.Rows.Item(648).Insert()

# Populate the newly inserted row 648 with its data
$ws.Cells.Item(648, 1).Value2 = 2
$ws.Cells.Item(648, 2).Value2 = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(648, 3).Value2 = "Coquimbo"
$ws.Cells.Item(648, 4).Value2 = 45218
$ws.Cells.Item(648, 5).Value2 = 4
$ws.Cells.Item(648, 6).Value2 = 100112043
$ws.Cells.Item(648, 7).Value2 = "Pepino dulce"
$ws.Cells.Item(648, 8).Value2 = "Sin especificar"
$ws.Cells.Item(648, 9).Value2 = "Segunda"
$ws.Cells.Item(648, 10).Value2 = 160
$ws.Cells.Item(648, 11).Value2 = 18000
$ws.Cells.Item(648, 12).Value2 = 19000
$ws.Cells.Item(648, 13).Value2 = 18500
$ws.Cells.Item(648, 14).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(648, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(648, 16).Value2 = 1028
$ws.Cells.Item(648, 17).Value2 = 18
$ws.Cells.Item(648, 18).Value2 = "Hortaliza"
